$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 536.3929000000001
$ws.Range("J17").Value = 552.55554
$ws.Range("L17").Value = 1657.66662
$ws.Range("N17").Value = -1993.66662
$ws.Range("H19").Value = 1941.9231
$ws.Range("I19").Value = 1066.5
$ws.Range("J19").Value = 2692.2856
$ws.Range("K19").Value = 1066.5
$ws.Range("L19").Value = 2692.2856
$ws.Range("M19").Value = -891.5
$ws.Range("N19").Value = -3042.2856
$ws.Range("H28").Value = 625462.6
$ws.Range("I28").Value = 769543.75
$ws.Range("K28").Value = 769543.75
$ws.Range("M28").Value = -769058.75
$ws.Range("H106").Value = 13657.129
$ws.Range("I106").Value = 14886.714
$ws.Range("K106").Value = 14886.714
$ws.Range("M106").Value = -14255.714
$ws.Range("H137").Value = 3000
$ws.Range("J137").Value = 5000
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -20100
$ws.Range("H138").Value = 3202.9622
$ws.Range("I138").Value = 2290.625
$ws.Range("J138").Value = 3597.4866
$ws.Range("K138").Value = 6871.875
$ws.Range("L138").Value = 10792.4598
$ws.Range("M138").Value = -1731.875
$ws.Range("N138").Value = -21072.4598

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2172.1365
$ws.Range("I2").Value = 1354.2354
$ws.Range("K2").Value = 1354.2354
$ws.Range("M2").Value = -1241.2354
$ws.Range("H32").Value = 1938.9323
$ws.Range("I32").Value = 896.1111
$ws.Range("K32").Value = 896.1111
$ws.Range("M32").Value = -609.1111
$ws.Range("H61").Value = 3818.3635
$ws.Range("J61").Value = 4772
$ws.Range("L61").Value = 4772
$ws.Range("N61").Value = -5196
$ws.Range("H69").Value = 307870
$ws.Range("J69").Value = 307870
$ws.Range("L69").Value = 307870
$ws.Range("N69").Value = -309368
$ws.Range("H72").Value = 307870
$ws.Range("J72").Value = 307870
$ws.Range("L72").Value = 923610
$ws.Range("N72").Value = -931098
$ws.Range("H110").Value = 1615
$ws.Range("I110").Value = 1615.0869
$ws.Range("J110").Value = 1613
$ws.Range("K110").Value = 1615.0869
$ws.Range("L110").Value = 1613
$ws.Range("M110").Value = 429.9131
$ws.Range("N110").Value = -5703
$ws.Range("H116").Value = 2172.1365
$ws.Range("I116").Value = 1354.2354
$ws.Range("K116").Value = 1354.2354
$ws.Range("M116").Value = 939.7646
$ws.Range("H124").Value = 31999
$ws.Range("J124").Value = 31999
$ws.Range("L124").Value = 31999
$ws.Range("N124").Value = -41819
$ws.Range("H132").Value = 7188.9585
$ws.Range("I132").Value = 7024.6045
$ws.Range("K132").Value = 21073.8135
$ws.Range("M132").Value = -18543.8135
$ws.Range("H136").Value = 3818.3635
$ws.Range("J136").Value = 4772
$ws.Range("L136").Value = 14316
$ws.Range("N136").Value = -19416

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2172.1365
$ws.Range("I3").Value = 1354.2354
$ws.Range("K3").Value = 1354.2354
$ws.Range("M3").Value = -1240.2354
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("K10").Value = 2000
$ws.Range("M10").Value = -1860
$ws.Range("H80").Value = 3041.25
$ws.Range("J80").Value = 5151
$ws.Range("L80").Value = 5151
$ws.Range("N80").Value = -7147
$ws.Range("H83").Value = 3041.25
$ws.Range("J83").Value = 5151
$ws.Range("L83").Value = 25755
$ws.Range("N83").Value = -35739
$ws.Range("H134").Value = 3408.383
$ws.Range("I134").Value = 3463.7
$ws.Range("J134").Value = 3092.2856
$ws.Range("K134").Value = 10391.1
$ws.Range("L134").Value = 9276.856800000001
$ws.Range("M134").Value = -7856.099999999999
$ws.Range("N134").Value = -14346.8568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6047.45
$ws.Range("I31").Value = 3742.7144
$ws.Range("J31").Value = 7288.4614
$ws.Range("K31").Value = 3742.7144
$ws.Range("L31").Value = 7288.4614
$ws.Range("M31").Value = -3447.7144
$ws.Range("N31").Value = -7878.4614
$ws.Range("H34").Value = 6047.45
$ws.Range("I34").Value = 3742.7144
$ws.Range("J34").Value = 7288.4614
$ws.Range("K34").Value = 3742.7144
$ws.Range("L34").Value = 7288.4614
$ws.Range("M34").Value = -3540.7144
$ws.Range("N34").Value = -7692.4614
$ws.Range("H41").Value = 31930.4
$ws.Range("J41").Value = 46999
$ws.Range("L41").Value = 46999
$ws.Range("N41").Value = -47855
$ws.Range("H105").Value = 1590
$ws.Range("I105").Value = 1663.4166
$ws.Range("K105").Value = 1663.4166
$ws.Range("M105").Value = 83.58339999999998
$ws.Range("H122").Value = 1945.5714
$ws.Range("I122").Value = 2001.25
$ws.Range("J122").Value = 1871.3334
$ws.Range("K122").Value = 6003.75
$ws.Range("L122").Value = 5614.0002
$ws.Range("M122").Value = -3553.75
$ws.Range("N122").Value = -10514.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35851196
$ws.Range("J4").Value = 41664.668
$ws.Range("L4").Value = 124994.004
$ws.Range("N4").Value = -125218.004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7031.76
$ws.Range("I132").Value = 9032.154
$ws.Range("K132").Value = 27096.462
$ws.Range("M132").Value = -24566.462

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 90000
$ws.Range("J36").Value = 90000
$ws.Range("L36").Value = 90000
$ws.Range("N36").Value = -91124
$ws.Range("H87").Value = 120000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 120000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H125").Value = 69355
$ws.Range("J125").Value = 69355
$ws.Range("L125").Value = 69355
$ws.Range("N125").Value = -79195

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 55000
$ws.Range("J110").Value = 55000
$ws.Range("L110").Value = 55000
$ws.Range("N110").Value = -63180
$ws.Range("H113").Value = 1189.3704
$ws.Range("I113").Value = 1365.7059
$ws.Range("K113").Value = 4097.1177
$ws.Range("M113").Value = -1927.1177
$ws.Range("H122").Value = 12503227
$ws.Range("I122").Value = 15154690
$ws.Range("J122").Value = 3470.2856
$ws.Range("K122").Value = 45464070
$ws.Range("L122").Value = 10410.8568
$ws.Range("M122").Value = -45461620
$ws.Range("N122").Value = -15310.8568
$ws.Range("H125").Value = 29607.5
$ws.Range("J125").Value = 29607.5
$ws.Range("L125").Value = 29607.5
$ws.Range("N125").Value = -39447.5
$ws.Range("H132").Value = 4208.3
$ws.Range("I132").Value = 3835.375
$ws.Range("J132").Value = 5700
$ws.Range("K132").Value = 11506.125
$ws.Range("L132").Value = 17100
$ws.Range("M132").Value = -8976.125
$ws.Range("N132").Value = -22160
